$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure textual price strings (e.g. "1.000", "0.4823") are not
# auto-converted to numbers by Excel, which would drop trailing
# zeros / change formatting. We mark the cell as Text format first
# only when the new value would otherwise parse as a plain number.

$ws.Range("D2").Value = "30.365.95"
$ws.Range("E2").Value = "  -1.29%  "

$ws.Range("D3").Value = "1.889.93"
$ws.Range("E3").Value = "  -1.46%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.95"
$ws.Range("E5").Value = "  -1.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4823"
$ws.Range("E7").Value = "  -1.98%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2900"
$ws.Range("E8").Value = "  -3.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06602"
$ws.Range("E9").Value = "  -2.70%  "

$ws.Range("D10").Value = "1.898.53"
$ws.Range("E10").Value = "  -1.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.91"
$ws.Range("E11").Value = "  -2.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07391"
$ws.Range("E12").Value = "  +0.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.175"
$ws.Range("E13").Value = "  -1.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.61"
$ws.Range("E14").Value = "  -1.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6622"
$ws.Range("E15").Value = "  -2.47%  "

$ws.Range("D16").Value = "30.367.70"

$ws.Range("E17").Value = "  -1.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007766"
$ws.Range("E18").Value = "  -3.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9996"
$ws.Range("E19").Value = "  -0.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.477"
$ws.Range("E20").Value = "  +0.40%  "

$ws.Range("D21").Value = "2.140.74"
$ws.Range("E21").Value = "  -1.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9990"
$ws.Range("E22").Value = "  -0.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "193.32"
$ws.Range("E23").Value = "  -3.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.189"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.434"
$ws.Range("E25").Value = "  -2.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.99"
$ws.Range("E26").Value = "  +1.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.29"
$ws.Range("E27").Value = "  -3.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.962"
$ws.Range("E28").Value = "  -0.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.449"
$ws.Range("E29").Value = "  -2.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.268"
$ws.Range("E30").Value = "  -2.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09157"
$ws.Range("E31").Value = "  -0.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.053"
$ws.Range("E32").Value = "  -0.95%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05090"
$ws.Range("E33").Value = "  -4.38%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7341"
$ws.Range("E34").Value = "  -1.98%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.143"
$ws.Range("E35").Value = "  +1.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.713"
$ws.Range("E36").Value = "  +0.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01814"
$ws.Range("E37").Value = "  -2.78%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.651"
$ws.Range("E38").Value = "  -2.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9200"
$ws.Range("E39").Value = "  -1.40%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.086"
$ws.Range("E40").Value = "  -0.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.65"
$ws.Range("E41").Value = "  -1.03%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4343"
$ws.Range("E42").Value = "  -3.84%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.893"
$ws.Range("E43").Value = "  -1.37%  "

$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.586"
$ws.Range("E45").Value = "  -2.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1333"
$ws.Range("E46").Value = "  -5.01%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "64.85"
$ws.Range("E47").Value = "  -10.75%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.548"
$ws.Range("E48").Value = "  +5.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.000"
$ws.Range("E49").Value = "  -1.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.20"
$ws.Range("E50").Value = "  -4.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05781"
$ws.Range("E51").Value = "  -2.39%  "
